# Apply the edits described by the commit:
# "Added all drop-down and non-dropdown classes for all sites."
#
# The underlying row/column data in Sheet1 stays the same conceptually,
# but three of the "class name" labels in column A change their letter
# casing (and the shared-string table gets reordered as a result):
#   A2: mdaTextHomepage -> mdaTextHomePage
#   A4: MdaTitle        -> mdaTitle
#   A8: pageTitlenewTab -> pageTitleNewTab
#
# The selected cell in the sheet also moves from B5 to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"

$ws.Range("A8").Select()
